# 平面关注对象.xlsx — "Merge loads in different storeys" edit
#
# Summary of the change (from the OOXML diff):
#   - Filter sheet gains a new row 21: Block="E-BFAS23-3", Filtering method="Ignore",
#     Note="液位控制器" (everything below row 20 shifts down by one, dimension
#     grows from D63 to D64).
#   - The Filter sheet becomes the active/selected tab (was Block before),
#     with its selection left on the newly inserted row (A21).
#   - The Block sheet's selection moves on to A118 (the sheet is no longer
#     the active tab).

$wb = $excel.ActiveWorkbook

$wsBlock  = $wb.Worksheets.Item("Block")
$wsFilter = $wb.Worksheets.Item("Filter")

# --- Filter sheet: insert the new row for E-BFAS23-3 --------------------
# Inserting a whole row at 21 pushes the former row 21 ("E-BFAS630" /
# "液位控制器" / ...) and everything after it down by one, which is exactly
# what the diff shows for rows 21-64.
$wsFilter.Rows.Item(21).Insert()

$wsFilter.Range("A21").Value = "E-BFAS23-3"
$wsFilter.Range("C21").Value = "Ignore"
$wsFilter.Range("D21").Value = "液位控制器"

# --- Selections / active tab --------------------------------------------
# Block keeps a selection, but is no longer the active sheet.
$wsBlock.Range("A118").Select()

# Filter becomes the active sheet, selection on the freshly inserted row.
$wsFilter.Activate()
$wsFilter.Range("A21").Select()
